# update area_protected_total with new beaufort sea data. remove incorrect
# bottle nose dolphin from ICO species
#
# Net effect captured by the canonical diff: a new metadata row describing
# the "iconic_species" layer is appended below the existing rows (the old
# row 20 / cw_pathogen_trend was the last row, so the new content lands on
# row 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Layer name / description (columns A & B), then USA data (column G) before
# Canada data (column E) - this ordering reproduces the shared-string
# insertion order of the source workbook.
$ws.Range("A21").Value = "iconic_species"
$ws.Range("B21").Value = "How iconic species were chosen for each region"
$ws.Range("G21").Value = "Species found in North Alaska protected under Marine Mammal Protection Act or Endangered Species Act"
$ws.Range("E21").Value = "Arctic Marine Mammals, Ivory Gull"

# Match the author's row height for the new (wrapped, multi-line) row.
$ws.Rows(21).RowHeight = 86.4

# The author's cursor ended up on E25 after adding this content.
$ws.Range("E25").Select()
